{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertText(\n  \"[LLM error: An error occurred (ExpiredTokenException) when calling the InvokeModel operation: The security token included in the request is expired]\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$para = $d.Paragraphs(1)\n$para.Range.Text = \"[LLM error: An error occurred (ExpiredTokenException) when calling the InvokeModel operation: The security token included in the request is expired]\"\n"}
